# Applies the automated crypto-price refresh (GitHub Actions run).
# Updates Price (D) and Volume(1h) (E) columns for each coin row, and
# fixes the Uniswap / BitcoinCash rows which swapped rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 20 and 21 swapped places (Uniswap <-> BitcoinCash) and got new values
$ws.Range("B20").Value2 = "BitcoinCash"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("B21").Value2 = "Uniswap"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"

# Price column (D) updates - cells whose new text looks numeric must be
# forced to Text format first so Excel keeps the exact original string
# (otherwise values like '0.5090' or '1.240' lose trailing zeros / become numbers).
$ws.Range("D2").Value2 = "25.914.73"
$ws.Range("D3").Value2 = "1.634.48"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.001"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "215.82"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.5090"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "1.001"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2581"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06345"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "19.53"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.07776"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "4.266"
$ws.Range("D13").Value2 = "1.635.31"
$ws.Range("D14").Value2 = "1.858.22"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.5519"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "63.88"
$ws.Range("D17").Value2 = "0.0₅7668"
$ws.Range("D18").Value2 = "25.920.58"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "1.002"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "195.10"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "4.422"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "9.902"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "6.051"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "1.001"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "1.915"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "142.19"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.1249"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "6.773"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "15.59"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.240"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.04891"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.247"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "3.199"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.544"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "2.367"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.8985"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.5531"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.542"
$ws.Range("D39").Value2 = "1.117.92"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.01556"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.9997"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "5.618"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "97.52"
$ws.Range("D46").Value2 = "1.768.87"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.4446"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.005"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "54.87"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.05135"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "7.592"

# Volume(1h) column (E) updates
$ws.Range("E2").Value2 = "  +0.10%  "
$ws.Range("E3").Value2 = "  +0.15%  "
$ws.Range("E4").Value2 = "  -0.23%  "
$ws.Range("E5").Value2 = "  +0.74%  "
$ws.Range("E6").Value2 = "  -0.42%  "
$ws.Range("E7").Value2 = "  -0.12%  "
$ws.Range("E8").Value2 = "  +1.32%  "
$ws.Range("E9").Value2 = "  +0.24%  "
$ws.Range("E10").Value2 = "  +0.66%  "
$ws.Range("E11").Value2 = "  +0.36%  "
$ws.Range("E12").Value2 = "  +0.14%  "
$ws.Range("E13").Value2 = "  -0.04%  "
$ws.Range("E14").Value2 = "  +0.02%  "
$ws.Range("E15").Value2 = "  +2.16%  "
$ws.Range("E16").Value2 = "  -0.06%  "
$ws.Range("E17").Value2 = "  -0.49%  "
$ws.Range("E18").Value2 = "  +0.13%  "
$ws.Range("E19").Value2 = "  -0.13%  "
$ws.Range("E20").Value2 = "  +0.52%  "
$ws.Range("E21").Value2 = "  +0.31%  "
$ws.Range("E22").Value2 = "  +0.21%  "
$ws.Range("E23").Value2 = "  +0.65%  "
$ws.Range("E24").Value2 = "  -0.27%  "
$ws.Range("E25").Value2 = "  +3.03%  "
$ws.Range("E26").Value2 = "  +1.04%  "
$ws.Range("E27").Value2 = "  +5.19%  "
$ws.Range("E28").Value2 = "  -0.45%  "
$ws.Range("E29").Value2 = "  +0.25%  "
$ws.Range("E30").Value2 = "  +0.42%  "
$ws.Range("E31").Value2 = "  +0.06%  "
$ws.Range("E32").Value2 = "  +0.39%  "
$ws.Range("E33").Value2 = "  +1.63%  "
$ws.Range("E34").Value2 = "  +1.33%  "
$ws.Range("E35").Value2 = "  +0.14%  "
$ws.Range("E36").Value2 = "  +1.45%  "
$ws.Range("E37").Value2 = "  +2.87%  "
$ws.Range("E38").Value2 = "  -1.31%  "
$ws.Range("E39").Value2 = "  -1.39%  "
$ws.Range("E40").Value2 = "  +0.80%  "
$ws.Range("E41").Value2 = "  -0.21%  "
$ws.Range("E42").Value2 = "  +3.29%  "
$ws.Range("E43").Value2 = "  -1.68%  "
$ws.Range("E44").Value2 = "  -1.13%  "
$ws.Range("E45").Value2 = "  -4.42%  "
$ws.Range("E46").Value2 = "  +0.03%  "
$ws.Range("E47").Value2 = "  -1.91%  "
$ws.Range("E48").Value2 = "  +0.17%  "
$ws.Range("E49").Value2 = "  +0.76%  "
$ws.Range("E50").Value2 = "  +1.65%  "
$ws.Range("E51").Value2 = "  +3.73%  "
